# Update "访问量"/count values in column F across the four sheets of the
# workbook to match the refreshed data snapshot generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 3214
$ws.Cells.Item(4, 6).Value = 1992
$ws.Cells.Item(5, 6).Value = 263
$ws.Cells.Item(7, 6).Value = 3071
$ws.Cells.Item(8, 6).Value = 609
$ws.Cells.Item(13, 6).Value = 153
$ws.Cells.Item(14, 6).Value = 153
$ws.Cells.Item(15, 6).Value = 10097
$ws.Cells.Item(17, 6).Value = 233
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(20, 6).Value = 7985
$ws.Cells.Item(21, 6).Value = 12589
$ws.Cells.Item(25, 6).Value = 267
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(29, 6).Value = 300
$ws.Cells.Item(30, 6).Value = 2813
$ws.Cells.Item(33, 6).Value = 7927
$ws.Cells.Item(34, 6).Value = 1466
$ws.Cells.Item(39, 6).Value = 1381
$ws.Cells.Item(42, 6).Value = 82

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(15, 6).Value = 14

# ---- Sheet "本地生活" ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 225
$ws.Cells.Item(5, 6).Value = 17

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 225
$ws.Cells.Item(5, 6).Value = 3214
$ws.Cells.Item(7, 6).Value = 1992
$ws.Cells.Item(9, 6).Value = 263
$ws.Cells.Item(10, 6).Value = 17
$ws.Cells.Item(11, 6).Value = 3071
$ws.Cells.Item(12, 6).Value = 609
$ws.Cells.Item(16, 6).Value = 153
$ws.Cells.Item(17, 6).Value = 153
$ws.Cells.Item(18, 6).Value = 10097
$ws.Cells.Item(19, 6).Value = 233
$ws.Cells.Item(20, 6).Value = 8
$ws.Cells.Item(22, 6).Value = 7986
$ws.Cells.Item(23, 6).Value = 12589
$ws.Cells.Item(27, 6).Value = 267
$ws.Cells.Item(31, 6).Value = 5
$ws.Cells.Item(32, 6).Value = 2813
$ws.Cells.Item(36, 6).Value = 7927
